$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 82.98768099999999
$ws.Range("H2").Value = 248.963043
$ws.Range("I2").Value = 0.4489504115427952
$ws.Range("J2").Value = 0.4489504115427952
$ws.Range("M2").Value = 560.2199806666666
$ws.Range("N2").Value = 1680.659942
$ws.Range("O2").Value = 0.6936344353529325
$ws.Range("P2").Value = 0.6936344353529326
$ws.Range("Q2").Value = 46491.35704539149
$ws.Range("R2").Value = 418422.2134085235
$ws.Range("S2").Value = 0.3114074652119534
$ws.Range("T2").Value = 0.3114074652119535
$ws.Range("G3").Value = 82.98768099999999
$ws.Range("H3").Value = 248.963043
$ws.Range("I3").Value = 0.4489504115427952
$ws.Range("J3").Value = 0.4489504115427952
$ws.Range("O3").Value = 0.06994956469466522
$ws.Range("P3").Value = 0.06994956469466522
$ws.Range("Q3").Value = 4688.420905364504
$ws.Range("R3").Value = 42195.78814828054
$ws.Range("S3").Value = 0.03140388585690933
$ws.Range("T3").Value = 0.03140388585690933
$ws.Range("G4").Value = 82.98768099999999
$ws.Range("H4").Value = 248.963043
$ws.Range("I4").Value = 0.4489504115427952
$ws.Range("J4").Value = 0.4489504115427952
$ws.Range("O4").Value = 0.2364159999524024
$ws.Range("P4").Value = 0.2364159999524024
$ws.Range("Q4").Value = 15845.95588804331
$ws.Range("R4").Value = 142613.6029923898
$ws.Range("S4").Value = 0.1061390604739325
$ws.Range("T4").Value = 0.1061390604739325
$ws.Range("G5").Value = 63.14058933333333
$ws.Range("I5").Value = 0.3415807409566563
$ws.Range("J5").Value = 0.3415807409566563
$ws.Range("M5").Value = 560.2199806666666
$ws.Range("N5").Value = 1680.659942
$ws.Range("O5").Value = 0.6936344353529325
$ws.Range("P5").Value = 0.6936344353529326
$ws.Range("Q5").Value = 35372.61973560193
$ws.Range("R5").Value = 318353.5776204174
$ws.Range("S5").Value = 0.2369321643809066
$ws.Range("T5").Value = 0.2369321643809066
$ws.Range("G6").Value = 63.14058933333333
$ws.Range("I6").Value = 0.3415807409566563
$ws.Range("J6").Value = 0.3415807409566563
$ws.Range("O6").Value = 0.06994956469466522
$ws.Range("P6").Value = 0.06994956469466522
$ws.Range("S6").Value = 0.02389342413799931
$ws.Range("T6").Value = 0.02389342413799931
$ws.Range("G7").Value = 63.14058933333333
$ws.Range("I7").Value = 0.3415807409566563
$ws.Range("J7").Value = 0.3415807409566563
$ws.Range("O7").Value = 0.2364159999524024
$ws.Range("P7").Value = 0.2364159999524024
$ws.Range("S7").Value = 0.08075515243775042
$ws.Range("T7").Value = 0.08075515243775042
$ws.Range("I8").Value = 0.2094688475005485
$ws.Range("J8").Value = 0.2094688475005485
$ws.Range("M8").Value = 560.2199806666666
$ws.Range("N8").Value = 1680.659942
$ws.Range("O8").Value = 0.6936344353529325
$ws.Range("P8").Value = 0.6936344353529326
$ws.Range("Q8").Value = 21691.68515865445
$ws.Range("R8").Value = 195225.1664278901
$ws.Range("S8").Value = 0.1452948057600724
$ws.Range("T8").Value = 0.1452948057600725
$ws.Range("I9").Value = 0.2094688475005485
$ws.Range("J9").Value = 0.2094688475005485
$ws.Range("O9").Value = 0.06994956469466522
$ws.Range("P9").Value = 0.06994956469466522
$ws.Range("S9").Value = 0.01465225469975658
$ws.Range("T9").Value = 0.01465225469975658
$ws.Range("I10").Value = 0.2094688475005485
$ws.Range("J10").Value = 0.2094688475005485
$ws.Range("O10").Value = 0.2364159999524024
$ws.Range("P10").Value = 0.2364159999524024
$ws.Range("Q10").Value = 7393.320135305331
$ws.Range("S10").Value = 0.04952178704071945
$ws.Range("T10").Value = 0.04952178704071945
